# Apply crypto price/volume updates (GitHub Actions scheduled refresh).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns D (Price) and E (Volume(1h)) hold plain text in the sheet (values like
# "1.00" / "10.00" or "  +0.01%  " must stay text, not be coerced to numbers), so
# force text storage with a temporary "@" number format, then restore the style.
$dataRange = $ws.Range("D2:E51")
$dataRange.NumberFormat = "@"

$ws.Range("D2").Value = "68.774.32"
$ws.Range("E2").Value = "  -0.32%  "
$ws.Range("D3").Value = "3.844.21"
$ws.Range("E3").Value = "  +2.40%  "
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "602.63"
$ws.Range("E5").Value = "  +0.08%  "
$ws.Range("D6").Value = "162.97"
$ws.Range("E6").Value = "  -2.78%  "
$ws.Range("D7").Value = "3.842.36"
$ws.Range("E7").Value = "  +2.44%  "
$ws.Range("E8").Value = "  -0.01%  "
$ws.Range("E9").Value = "  -1.66%  "
$ws.Range("E10").Value = "  -0.90%  "
$ws.Range("D11").Value = "6.30"
$ws.Range("E11").Value = "  -2.58%  "
$ws.Range("E12").Value = "  -0.16%  "
$ws.Range("D13").Value = "36.90"
$ws.Range("E13").Value = "  -2.65%  "
$ws.Range("E14").Value = "  -2.00%  "
$ws.Range("D15").Value = "4.482.76"
$ws.Range("E15").Value = "  +2.31%  "
$ws.Range("D16").Value = "3.847.53"
$ws.Range("E16").Value = "  +2.50%  "
$ws.Range("D17").Value = "68.970.92"
$ws.Range("E17").Value = "  -0.09%  "
$ws.Range("E18").Value = "  +2.82%  "
$ws.Range("B19").Value = "TRON"
$ws.Range("C19").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D19").Value = "0.113"
$ws.Range("E19").Value = "  -0.18%  "
$ws.Range("B20").Value = "Uniswap"
$ws.Range("C20").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D20").Value = "11.38"
$ws.Range("E20").Value = "  +4.53%  "
$ws.Range("D21").Value = "17.13"
$ws.Range("E21").Value = "  -0.85%  "
$ws.Range("D22").Value = "484.62"
$ws.Range("D23").Value = "0.719"
$ws.Range("E23").Value = "  -0.95%  "
$ws.Range("E24").Value = "  +3.09%  "
$ws.Range("D25").Value = "84.03"
$ws.Range("E25").Value = "  -0.92%  "
$ws.Range("E26").Value = "  -2.18%  "
$ws.Range("D27").Value = "12.09"
$ws.Range("E27").Value = "  -1.95%  "
$ws.Range("D28").Value = "10.00"
$ws.Range("D29").Value = "0.999"
$ws.Range("E29").Value = "  -0.11%  "
$ws.Range("D30").Value = "2.97"
$ws.Range("E30").Value = "  -0.94%  "
$ws.Range("D31").Value = "7.93"
$ws.Range("E31").Value = "  -1.01%  "
$ws.Range("D32").Value = "3.998.61"
$ws.Range("E32").Value = "  +2.56%  "
$ws.Range("E33").Value = "  -3.93%  "
$ws.Range("D34").Value = "32.13"
$ws.Range("E34").Value = "  +1.48%  "
$ws.Range("D35").Value = "3.791.63"
$ws.Range("E35").Value = "  +2.81%  "
$ws.Range("E36").Value = "  -1.83%  "
$ws.Range("E37").Value = "  +1.03%  "
$ws.Range("D38").Value = "0.140"
$ws.Range("E38").Value = "  +4.65%  "
$ws.Range("D39").Value = "5.87"
$ws.Range("E39").Value = "  -0.46%  "
$ws.Range("D40").Value = "1.00"
$ws.Range("E40").Value = "  +0.03%  "
$ws.Range("E41").Value = "  -2.07%  "
$ws.Range("D42").Value = "441.86"
$ws.Range("E42").Value = "  +2.90%  "
$ws.Range("D43").Value = "2.98"
$ws.Range("E43").Value = "  +1.03%  "
$ws.Range("D44").Value = "48.53"
$ws.Range("E44").Value = "  -0.23%  "
$ws.Range("E45").Value = "  -1.54%  "
$ws.Range("E46").Value = "  +0.00%  "
$ws.Range("D47").Value = "8.39"
$ws.Range("E47").Value = "  -1.15%  "
$ws.Range("D48").Value = "27.34"
$ws.Range("E48").Value = "  +16.95%  "
$ws.Range("D49").Value = "2.833.18"
$ws.Range("E49").Value = "  +1.48%  "
$ws.Range("D50").Value = "142.34"
$ws.Range("E50").Value = "  +0.55%  "
$ws.Range("D51").Value = "0.0355"
$ws.Range("E51").Value = "  +0.86%  "

$dataRange.Style = "Normal"

